# Error Calculations and Plots
#
# Two data rows are removed entirely (RM 232 and SC 92), shifting all
# subsequent rows up, and a number of individual cells are updated to
# new computed/imputed values (or cleared back to missing).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the higher-numbered row first so the lower row's index doesn't move.
$ws.Rows("28:28").Delete()
$ws.Rows("26:26").Delete()

# Apply the remaining per-cell value changes (using final, post-delete row numbers).
$ws.Range("E2").Value = -7.2
$ws.Range("F3").ClearContents()
$ws.Range("F4").Value = 17.97
$ws.Range("E6").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("F9").ClearContents()
$ws.Range("E12").Value = -5.3
$ws.Range("E14").ClearContents()
$ws.Range("F15").Value = 16.2
$ws.Range("F18").Value = 18.35
$ws.Range("F19").ClearContents()
$ws.Range("E20").Value = -7.2
$ws.Range("E21").Value = -8.699999999999999
$ws.Range("F22").ClearContents()
$ws.Range("E23").ClearContents()
$ws.Range("F23").Value = 16.48
$ws.Range("E24").ClearContents()
$ws.Range("F25").Value = 16.6
$ws.Range("D26").Value = -13.8
$ws.Range("D27").ClearContents()
$ws.Range("F27").ClearContents()
$ws.Range("D28").ClearContents()
$ws.Range("D29").Value = -13
$ws.Range("D30").Value = -13.6
$ws.Range("D31").ClearContents()
$ws.Range("E31").Value = -8.1
$ws.Range("D32").ClearContents()
$ws.Range("E33").Value = -10.7
